$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K. This shifts the existing K:O (Created Date,
# Created By, Modified Date, Modified By, Status) one column to the right
# (L:P), and automatically extends the A1:O1 title merge to A1:P1 and the
# used-range dimension from O3 to P3.
$ws.Columns("K").Insert()

# New "Image" header + uploaded image path value for the vehicle spec row.
$ws.Range("K2").Value = "Image"
$ws.Range("K3").Value = "C:\Users\pc\Pictures\FMS\FlowChart\CAF_Page_Edited.png"

# Make sure the new data cell has the same bordered look as its row neighbours.
$ws.Range("C3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Series value updated from "GLS 2.0 A/T" to "PX 2.0 A/T".
$ws.Range("C3").Value = "PX 2.0 A/T"

# Selection moved to E3.
$ws.Range("E3").Select()
